$wb = $excel.ActiveWorkbook

# Locate the current last sheet (StoreLocation) so the new sheet is added after it.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# Add the new "ContactSales" worksheet right after the last existing sheet.
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "ContactSales"

# Populate row 2 with the contact-sales form field values.
$ws.Range("A2").Value = "us"
$ws.Range("B2").Value = "10-19"
$ws.Range("C2").Value = "Advanced_Communications"
$ws.Range("D2").Value = "technology"
$ws.Range("E2").Value = "email"
$ws.Range("F2").Value = "NY"
$ws.Range("G2").Value = "Submit"

# Copy the existing text formatting (Menlo font, text number format) from the
# StoreLocation sheet's B2 cell so the new cells share the same style.
$srcCell = $lastSheet.Range("B2")
$srcCell.Copy()
$dst = $ws.Range("A2:G2")
$dst.PasteSpecial(-4122)  # xlPasteFormats

# Match the column widths from the authored workbook.
$ws.Columns.Item(3).ColumnWidth = 28.5
$ws.Columns.Item(4).ColumnWidth = 12.666666666666666

# Select the new row, mirroring the saved selection state.
$ws.Range("A2:G2").Select() | Out-Null

# Make the new sheet the active tab.
$ws.Activate() | Out-Null
